# Update odds values on Sheet1 to reflect the latest Betfair Back/Lay
# snapshot for 2025-12-25 ("Atualizando o arquivo XLSX").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Al Jabalain vs Al-Ula FC)
$ws.Range("V2").Value = 1.78

# Row 3 (Al-Feiha vs Al-Hazm (KSA))
$ws.Range("Q3").Value = 1.87

# Row 5 (Al-Wahda (KSA) vs Al Jubail)
$ws.Range("G5").Value = 1.69
$ws.Range("S5").Value = 3.05
$ws.Range("W5").Value = 2.44
$ws.Range("X5").Value = 990
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 980
$ws.Range("AJ5").Value = 17.5

# Row 6 (Al Riyadh SC vs Al-Ettifaq)
$ws.Range("I6").Value = 2.58
$ws.Range("J6").Value = 3.3
$ws.Range("Q6").Value = 1.91
$ws.Range("V6").Value = 1.63
$ws.Range("AE6").Value = 36

# Row 7 (NEOM Sports Club vs Al Najma Club)
$ws.Range("F7").Value = 1.42
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 4.5
$ws.Range("O7").Value = 1.22
$ws.Range("R7").Value = 1.48
$ws.Range("S7").Value = 2.64
$ws.Range("T7").Value = 1.89
$ws.Range("U7").Value = 1.9
$ws.Range("X7").Value = 990
$ws.Range("Z7").Value = 80
$ws.Range("AA7").Value = 320
$ws.Range("AB7").Value = 10.5
$ws.Range("AC7").Value = 1000
$ws.Range("AE7").Value = 150
$ws.Range("AF7").Value = 9.6
$ws.Range("AG7").Value = 11
$ws.Range("AI7").Value = 130
$ws.Range("AJ7").Value = 1000
$ws.Range("AK7").Value = 1000
$ws.Range("AM7").Value = 160
$ws.Range("AN7").Value = 6.4

# Row 8 (Belouizdad vs ES Setif)
$ws.Range("F8").Value = 1.63
$ws.Range("G8").Value = 1.8
$ws.Range("H8").Value = 5.8
$ws.Range("J8").Value = 3
$ws.Range("K8").Value = 4.3
$ws.Range("T8").Value = 2.3
$ws.Range("V8").Value = 1.12
$ws.Range("W8").Value = 2.24

# Row 9 (Olancho vs Platense FC)
$ws.Range("J9").Value = 1.03
$ws.Range("L9").Value = 1.01
$ws.Range("M9").Value = 1.01
$ws.Range("N9").Value = 1.26
$ws.Range("O9").Value = 1.28
$ws.Range("P9").Value = 1.26
$ws.Range("Q9").Value = 1.28
$ws.Range("R9").Value = 1.18
$ws.Range("S9").Value = 1.28
$ws.Range("T9").Value = 1.01
$ws.Range("U9").Value = 1.01
$ws.Range("V9").Value = 1.01
$ws.Range("W9").Value = 1.01
$ws.Range("X9").Value = 1000
$ws.Range("Y9").Value = 1000
$ws.Range("Z9").Value = 1000
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 1000
$ws.Range("AC9").Value = 1000
$ws.Range("AD9").Value = 1000
$ws.Range("AE9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AG9").Value = 1000
$ws.Range("AH9").Value = 1000
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 1000
$ws.Range("AK9").Value = 1000
$ws.Range("AL9").Value = 1000
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 1000
$ws.Range("AO9").Value = 1000
